# Applies corrected Diebold-Mariano statistics (DM_Stat, P_Value) values
# to columns C and D for rows 2-11 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  C = -1.157645537839379;  D = 0.2594229912477282 },
    @{ Row = 3;  C = 0.344936218188253;   D = 0.7334182632280042 },
    @{ Row = 4;  C = 0.4346158584302625;  D = 0.6680745087073683 },
    @{ Row = 5;  C = 0.9180471910989026;  D = 0.3685499601385649 },
    @{ Row = 6;  C = 1.32462201449687;    D = 0.1988926119763745 },
    @{ Row = 7;  C = 1.236481384496879;   D = 0.2293181007223457 },
    @{ Row = 8;  C = 1.591433687749416;   D = 0.1257800612360589 },
    @{ Row = 9;  C = 0.02221119352114896; D = 0.982479715355711 },
    @{ Row = 10; C = 0.7360613980450912;  D = 0.4694694395843062 },
    @{ Row = 11; C = 0.6578989152028353;  D = 0.5174286420002172 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
